$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Commitment date values for rows 2-7 (2023-01-20 -> serial 44946)
$commitDate = Get-Date -Year 2023 -Month 1 -Day 20 -Hour 0 -Minute 0 -Second 0

# Apply the date format to J2 first, then propagate the identical style to
# J3:J7 via copy/paste-format so all six cells share a single cellXf entry
# (matches the source workbook, which reuses one new style for the column).
$ws.Range("J2").NumberFormat = "mm-dd-yy"
$ws.Range("J2").Copy()
$ws.Range("J3:J7").PasteSpecial(-4122)
$ws.Range("J2:J7").Value = $commitDate

# Autofit the new column to the date values, then add the header.
$ws.Columns.Item(10).AutoFit()
$ws.Range("J1").Value = "Commitment Date"

# Match the new selection left behind by the edit.
$ws.Range("J3:J7").Select()
